$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.026.78'
$ws.Range("E2").Value = '  +0.61%  '

$ws.Range("D3").Value = '1.920.53'
$ws.Range("E3").Value = '  +1.64%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.26'
$ws.Range("E5").Value = '  +0.62%  '

$ws.Range("E6").Value = '  +0.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4593'
$ws.Range("E7").Value = '  +0.23%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3813'
$ws.Range("E8").Value = '  +0.17%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07747'
$ws.Range("E9").Value = '  +0.37%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9783'
$ws.Range("E10").Value = '  +1.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.60'
$ws.Range("E11").Value = '  +2.55%  '

$ws.Range("D12").Value = '1.915.90'
$ws.Range("E12").Value = '  +1.93%  '

$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.684'
$ws.Range("E13").Value = '  +0.48%  '

$ws.Range("B14").Value = 'Chainlink'
$ws.Range("C14").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.960'
$ws.Range("E14").Value = '  +0.21%  '

$ws.Range("E15").Value = '  -0.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.006'
$ws.Range("E16").Value = '  +0.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '84.23'
$ws.Range("E17").Value = '  +1.29%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009522'
$ws.Range("E18").Value = '  +0.27%  '

$ws.Range("E19").Value = '  +0.68%  '

$ws.Range("E20").Value = '  +0.14%  '

$ws.Range("D21").Value = '29.034.91'
$ws.Range("E21").Value = '  +0.77%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.338'
$ws.Range("E22").Value = '  +0.65%  '

$ws.Range("E23").Value = '  +1.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.075'
$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.45'
$ws.Range("E25").Value = '  +0.74%  '

$ws.Range("E26").Value = '  +0.17%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.650'
$ws.Range("E27").Value = '  +1.00%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '117.68'
$ws.Range("E28").Value = '  +0.50%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.836'
$ws.Range("E29").Value = '  +1.53%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.09331'
$ws.Range("E30").Value = '  +0.67%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.8592'
$ws.Range("E31").Value = '  +0.90%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.098'
$ws.Range("E32").Value = '  +0.50%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.241'
$ws.Range("E33").Value = '  +1.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.018'
$ws.Range("E34").Value = '  +0.15%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.158'
$ws.Range("E35").Value = '  +1.51%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05679'
$ws.Range("E36").Value = '  +0.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.004'
$ws.Range("E37").Value = '  +0.16%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02044'
$ws.Range("E38").Value = '  +0.61%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.100'
$ws.Range("E39").Value = '  +15.17%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.437'
$ws.Range("E40").Value = '  +0.46%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5497'
$ws.Range("E41").Value = '  +0.22%  '

$ws.Range("E42").Value = '  +0.30%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.364'
$ws.Range("E43").Value = '  +2.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.192'
$ws.Range("E44").Value = '  +5.97%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.000002759'
$ws.Range("E45").Value = '  -4.86%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5183'
$ws.Range("E46").Value = '  +0.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '11.25'
$ws.Range("E47").Value = '  +0.27%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.06909'
$ws.Range("E48").Value = '  +1.57%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '110.50'
$ws.Range("E49").Value = '  -1.01%  '

$ws.Range("E50").Value = '  -0.82%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.004'
$ws.Range("E51").Value = '  +0.20%  '
